$wb = $excel.ActiveWorkbook

# Updated "want to go" counts (column F) for rows 2-5, applied to both the
# "展览" and "全部类型" sheets, which carry duplicate data.
$sheetNames = @("展览", "全部类型")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 61
    $ws.Range("F3").Value = 416
    $ws.Range("F4").Value = 27
    $ws.Range("F5").Value = 125
}
